# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
#
# Updates the IG StructureDefinition export (child-organization-hierarchy-level-code)
# from version 5.0.0 -> 6.0.0:
#   - Metadata sheet: Version, Date, Publisher, and replaces the duplicated
#     "Contact / No display for ContactDetail" rows with a single
#     "Publisher" value row plus a new "Jurisdiction" row.
#   - Elements sheet: the root Extension row's Type(s)/Short columns now
#     describe the concrete extension instead of the generic "Extension".

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value fills in (was blank) as "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Drop the duplicated second "Contact" row (old rows 10 & 11 both held
# Contact / "No display for ContactDetail"); delete one of them.
$meta.Rows.Item(11).Delete()

# The remaining row becomes a new "Jurisdiction" / "United States of America" pair.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet --------------------------------------------------------

# Root Extension row: Type(s) and Short description become specific to
# this extension rather than the generic Extension placeholders.
$elements.Cells.Item(2, 11).Value = "Child Organization Hierarchy Level Code"
$elements.Cells.Item(2, 12).Value = "Numeric level of the child practitioner within the organinzational hierarchy"

# The "Short" column widened (bestFit) to accommodate the new, longer text.
$elements.Columns.Item(11).ColumnWidth = 37.25
